# RPA datasets push 2023-11-11
#
# Two IPO listings that previously had no confirmed offering price ("-")
# now have their final confirmed price ("확정공모가"):
#   row 11 - 그린리소스        -> 17000
#   row 12 - 한선엔지니어링    -> 7000
#
# The column stores these prices as plain text (shared strings), the same
# way every other row in the column already does (e.g. D13 = "34700"), so
# a leading apostrophe forces the literal text entry instead of letting
# Excel convert it to a number. The cell style is then reset to "Normal"
# so the quote-prefix formatting doesn't linger on the cell (matching the
# rest of the column, which carries no explicit style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D11").Value = "'17000"
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").Value = "'7000"
$ws.Range("D12").Style = "Normal"
